$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Fix the horizontal score algorithm's header layout on Sheet1 ---

# Row 2: the middle header cells (HEADER2/HEADER3) are no longer needed
# now that a dedicated sub-header row is being inserted below.
$ws1.Range("C2").Value = ""
$ws1.Range("D2").Value = ""

# Row 3: drop the old SH2/SH4 sub-header labels; SH1/SH3/SH5 remain.
$ws1.Range("C3").Value = ""
$ws1.Range("E3").Value = ""

# Insert a new row at 4 for the expanded sub-header set (shifts the
# existing data rows down by one).
$ws1.Rows("4").Insert()

# Populate the new sub-header row.
$ws1.Range("B4").Value = "SH21"
$ws1.Range("C4").Value = "SH22"
$ws1.Range("D4").Value = "SH23"
$ws1.Range("E4").Value = "SH24"
$ws1.Range("F4").Value = "SH25"

# Give the whole header block (now 3 rows tall) a thin box border on
# every cell, on top of the existing bold font + fill.
$ws1.Range("B2:F4").Borders.LineStyle = 1

# Sheet1 becomes the active sheet/tab, with D4 selected.
$ws1.Activate()
[void]$ws1.Range("D4").Select()
